# Apply weekly update: insert two new rows (new week of data) above the
# previous row 18, pushing the existing rows 18-23 down to 20-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 18 (formatting is inherited from row 17,
# which already carries the date style used by column D).
$ws.Rows("18:19").Insert()

# --- New row 18: Perejil, Primera, week of 2022-10-05 (serial 44839) ---
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C18").Value = 'Ñuble'
$ws.Range("D18").Value = 44839
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112044
$ws.Range("G18").Value = 'Perejil'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 240
$ws.Range("K18").Value = 700
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = 750
$ws.Range("N18").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O18").Value = 'Región del Maule'
$ws.Range("P18").Value = 750
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 'Hortaliza'

# --- New row 19: Perejil, Segunda, week of 2022-10-05 (serial 44839) ---
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C19").Value = 'Ñuble'
$ws.Range("D19").Value = 44839
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112044
$ws.Range("G19").Value = 'Perejil'
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 600
$ws.Range("M19").Value = 600
$ws.Range("N19").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O19").Value = 'Región del Maule'
$ws.Range("P19").Value = 600
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 'Hortaliza'
